# Scheduled market-data refresh for the Leve profit-calculator workbook.
# Updates currentAveragePrice* / LeveProfit* columns (H,I,J,K,L,M,N) on each
# job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly polled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 629.8570999999999
$ws.Range("I2").Value = 401.66666
$ws.Range("K2").Value = 401.66666
$ws.Range("M2").Value = -288.66666
# Row 33
$ws.Range("H33").Value = 1684093.8
$ws.Range("I33").Value = 2469624
$ws.Range("J33").Value = 814.7143
$ws.Range("K33").Value = 2469624
$ws.Range("L33").Value = 814.7143
$ws.Range("M33").Value = -2469395
$ws.Range("N33").Value = -1272.7143
# Row 59
$ws.Range("H59").Value = 320
$ws.Range("J59").Value = 320
$ws.Range("L59").Value = 960
$ws.Range("N59").Value = -2074
# Row 70
$ws.Range("H70").Value = 4993.9443
$ws.Range("I70").Value = 927
$ws.Range("J70").Value = 7582
$ws.Range("K70").Value = 2781
$ws.Range("L70").Value = 22746
$ws.Range("M70").Value = -2511
$ws.Range("N70").Value = -23286
# Row 73
$ws.Range("H73").Value = 4993.9443
$ws.Range("I73").Value = 927
$ws.Range("J73").Value = 7582
$ws.Range("K73").Value = 2781
$ws.Range("L73").Value = 22746
$ws.Range("M73").Value = -1845
$ws.Range("N73").Value = -24618
# Row 80
$ws.Range("H80").Value = 3556.2104
$ws.Range("I80").Value = 471
$ws.Range("J80").Value = 6332.9
$ws.Range("K80").Value = 1413
$ws.Range("L80").Value = 18998.7
$ws.Range("M80").Value = -415
$ws.Range("N80").Value = -20994.7
# Row 83
$ws.Range("H83").Value = 3556.2104
$ws.Range("I83").Value = 471
$ws.Range("J83").Value = 6332.9
$ws.Range("K83").Value = 4239
$ws.Range("L83").Value = 56996.1
$ws.Range("M83").Value = 753
$ws.Range("N83").Value = -66980.10000000001
# Row 98
$ws.Range("H98").Value = 1317.2858
$ws.Range("I98").Value = 1321.4117
$ws.Range("J98").Value = 1299.75
$ws.Range("K98").Value = 1321.4117
$ws.Range("L98").Value = 1299.75
$ws.Range("M98").Value = 176.5882999999999
$ws.Range("N98").Value = -4295.75
# Row 100
$ws.Range("H100").Value = 867.6667
$ws.Range("I100").Value = 954.3
$ws.Range("K100").Value = 954.3
$ws.Range("M100").Value = -413.3
# Row 103
$ws.Range("H103").Value = 778.1053000000001
$ws.Range("J103").Value = 802.1429000000001
$ws.Range("L103").Value = 2406.4287
$ws.Range("N103").Value = -3578.4287
# Row 106
$ws.Range("H106").Value = 31251830
$ws.Range("I106").Value = 33335152
$ws.Range("K106").Value = 33335152
$ws.Range("M106").Value = -33334521
# Row 122
$ws.Range("H122").Value = 1317.2858
$ws.Range("I122").Value = 1321.4117
$ws.Range("J122").Value = 1299.75
$ws.Range("K122").Value = 3964.2351
$ws.Range("L122").Value = 3899.25
$ws.Range("M122").Value = -1514.2351
$ws.Range("N122").Value = -8799.25
# Row 132
$ws.Range("H132").Value = 5566.227
$ws.Range("I132").Value = 6370.0557
$ws.Range("K132").Value = 19110.1671
$ws.Range("M132").Value = -16580.1671
# Row 137
$ws.Range("H137").Value = 65930.96000000001
$ws.Range("I137").Value = 128885.93
$ws.Range("J137").Value = 2976
$ws.Range("K137").Value = 386657.79
$ws.Range("L137").Value = 8928
$ws.Range("M137").Value = -384107.79
$ws.Range("N137").Value = -14028
# Row 138
$ws.Range("H138").Value = 3314.6875
$ws.Range("I138").Value = 2506.7334
$ws.Range("J138").Value = 3562.0205
$ws.Range("K138").Value = 7520.2002
$ws.Range("L138").Value = 10686.0615
$ws.Range("M138").Value = -2380.2002
$ws.Range("N138").Value = -20966.0615
# Row 141
$ws.Range("H141").Value = 61583.57
$ws.Range("I141").Value = 40666.668
$ws.Range("J141").Value = 77271.25
$ws.Range("K141").Value = 122000.004
$ws.Range("L141").Value = 231813.75
$ws.Range("M141").Value = -116820.004
$ws.Range("N141").Value = -242173.75
$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 4066.25
$ws.Range("I30").Value = 7500
$ws.Range("J30").Value = 2921.6667
$ws.Range("K30").Value = 7500
$ws.Range("L30").Value = 2921.6667
$ws.Range("M30").Value = -7350
$ws.Range("N30").Value = -3221.6667
# Row 32
$ws.Range("H32").Value = 6810.513
$ws.Range("I32").Value = 5688.357
$ws.Range("J32").Value = 19902.334
$ws.Range("K32").Value = 5688.357
$ws.Range("L32").Value = 19902.334
$ws.Range("M32").Value = -5401.357
$ws.Range("N32").Value = -20476.334
# Row 34
$ws.Range("H34").Value = 15000
$ws.Range("J34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15542
# Row 61
$ws.Range("H61").Value = 4780.4375
$ws.Range("I61").Value = 4731.3076
$ws.Range("K61").Value = 4731.3076
$ws.Range("M61").Value = -4519.3076
# Row 74
$ws.Range("H74").Value = 93509.8
$ws.Range("I74").Value = 1842.2
$ws.Range("K74").Value = 1842.2
$ws.Range("M74").Value = -968.2
# Row 77
$ws.Range("H77").Value = 93509.8
$ws.Range("I77").Value = 1842.2
$ws.Range("K77").Value = 9211
$ws.Range("M77").Value = -4843
# Row 110
$ws.Range("H110").Value = 897097.9
$ws.Range("I110").Value = 1208560.2
$ws.Range("J110").Value = 1643.375
$ws.Range("K110").Value = 1208560.2
$ws.Range("L110").Value = 1643.375
$ws.Range("M110").Value = -1206515.2
$ws.Range("N110").Value = -5733.375
# Row 122
$ws.Range("H122").Value = 2451925.2
$ws.Range("I122").Value = 2633157.2
$ws.Range("J122").Value = 2089461.2
$ws.Range("K122").Value = 7899471.600000001
$ws.Range("L122").Value = 6268383.6
$ws.Range("M122").Value = -7897021.600000001
$ws.Range("N122").Value = -6273283.6
# Row 132
$ws.Range("H132").Value = 3078
$ws.Range("I132").Value = 1944
$ws.Range("K132").Value = 5832
$ws.Range("M132").Value = -3302
# Row 133
$ws.Range("H133").Value = 98000
$ws.Range("J133").Value = 98000
$ws.Range("L133").Value = 98000
$ws.Range("N133").Value = -103060
# Row 136
$ws.Range("H136").Value = 4780.4375
$ws.Range("I136").Value = 4731.3076
$ws.Range("K136").Value = 14193.9228
$ws.Range("M136").Value = -11643.9228
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 5954343
$ws.Range("I99").Value = 7520117.5
$ws.Range("J99").Value = 4399.4
$ws.Range("K99").Value = 7520117.5
$ws.Range("L99").Value = 4399.4
$ws.Range("M99").Value = -7518619.5
$ws.Range("N99").Value = -7395.4
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20807.371
$ws.Range("I31").Value = 3085.3
$ws.Range("J31").Value = 31232.117
$ws.Range("K31").Value = 3085.3
$ws.Range("L31").Value = 31232.117
$ws.Range("M31").Value = -2790.3
$ws.Range("N31").Value = -31822.117
# Row 34
$ws.Range("H34").Value = 20807.371
$ws.Range("I34").Value = 3085.3
$ws.Range("J34").Value = 31232.117
$ws.Range("K34").Value = 3085.3
$ws.Range("L34").Value = 31232.117
$ws.Range("M34").Value = -2883.3
$ws.Range("N34").Value = -31636.117
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 48039.9
$ws.Range("J37").Value = 48039.9
$ws.Range("L37").Value = 144119.7
$ws.Range("N37").Value = -144343.7
# Row 46
$ws.Range("H46").Value = 130372.54
$ws.Range("I46").Value = 333489.4
$ws.Range("J46").Value = 3424.5
$ws.Range("K46").Value = 1000468.2
$ws.Range("L46").Value = 10273.5
$ws.Range("M46").Value = -1000377.2
$ws.Range("N46").Value = -10455.5
# Row 86
$ws.Range("H86").Value = 967.3333
$ws.Range("I86").Value = 999
$ws.Range("J86").Value = 951.5
$ws.Range("K86").Value = 2997
$ws.Range("L86").Value = 2854.5
$ws.Range("M86").Value = -1811
$ws.Range("N86").Value = -5226.5
# Row 89
$ws.Range("H89").Value = 967.3333
$ws.Range("I89").Value = 999
$ws.Range("J89").Value = 951.5
$ws.Range("K89").Value = 8991
$ws.Range("L89").Value = 8563.5
$ws.Range("M89").Value = -3063
$ws.Range("N89").Value = -20419.5
# Row 137
$ws.Range("H137").Value = 8194.25
$ws.Range("I137").Value = 7941.6665
$ws.Range("J137").Value = 8345.799999999999
$ws.Range("K137").Value = 23824.9995
$ws.Range("L137").Value = 25037.4
$ws.Range("M137").Value = -18724.9995
$ws.Range("N137").Value = -35237.39999999999
$ws = $wb.Worksheets.Item("GSM")
# Row 28
$ws.Range("H28").Value = 9750
$ws.Range("J28").Value = 9750
$ws.Range("L28").Value = 9750
$ws.Range("N28").Value = -10134
# Row 117
$ws.Range("H117").Value = 49988
$ws.Range("J117").Value = 49988
$ws.Range("L117").Value = 49988
$ws.Range("N117").Value = -56872
# Row 122
$ws.Range("H122").Value = 332530.8
$ws.Range("I122").Value = 389101.4
$ws.Range("K122").Value = 1167304.2
$ws.Range("M122").Value = -1164854.2
# Row 132
$ws.Range("H132").Value = 3451.6553
$ws.Range("I132").Value = 3328.75
$ws.Range("J132").Value = 3724.7778
$ws.Range("K132").Value = 9986.25
$ws.Range("L132").Value = 11174.3334
$ws.Range("M132").Value = -7456.25
$ws.Range("N132").Value = -16234.3334
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4630560.5
$ws.Range("I61").Value = 5291873
$ws.Range("K61").Value = 5291873
$ws.Range("M61").Value = -5291671
# Row 68
$ws.Range("H68").Value = 4049.8
$ws.Range("J68").Value = 4083
$ws.Range("L68").Value = 4083
$ws.Range("N68").Value = -5581
# Row 71
$ws.Range("H71").Value = 4049.8
$ws.Range("J71").Value = 4083
$ws.Range("L71").Value = 20415
$ws.Range("N71").Value = -27903
# Row 113
$ws.Range("H113").Value = 4630560.5
$ws.Range("I113").Value = 5291873
$ws.Range("K113").Value = 5291873
$ws.Range("M113").Value = -5289703
# Row 122
$ws.Range("H122").Value = 5853.8335
$ws.Range("I122").Value = 3865.3333
$ws.Range("K122").Value = 11595.9999
$ws.Range("M122").Value = -9145.999899999999
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2947.8438
$ws.Range("I122").Value = 1526.04
$ws.Range("K122").Value = 4578.12
$ws.Range("M122").Value = -2128.12
# Row 126
$ws.Range("H126").Value = 5031.5454
$ws.Range("I126").Value = 3316.6667
$ws.Range("J126").Value = 7089.4
$ws.Range("K126").Value = 9950.000100000001
$ws.Range("L126").Value = 21268.2
$ws.Range("M126").Value = -7480.000100000001
$ws.Range("N126").Value = -26208.2
# Row 136
$ws.Range("H136").Value = 1172.2325
$ws.Range("I136").Value = 1065.7646
$ws.Range("J136").Value = 1574.4445
$ws.Range("K136").Value = 3197.2938
$ws.Range("L136").Value = 4723.333500000001
$ws.Range("M136").Value = -647.2937999999999
$ws.Range("N136").Value = -9823.333500000001
